# UI update David ILunga
# Registre Canabis.xlsx - update "Menu déroulant:" placeholder value in F2
# to the actual selected dropdown option, and refresh the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "micro dissection"

# Reflect the new active cell/selection left after the edit.
$ws.Range("F2").Select() | Out-Null
